$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: "Eastbourne, Newhaven & Seaford" (Anthony Sheehy) ---
# Mark as added to map, attach the project image as a real hyperlink, note "-"
$ws.Range("C28").Value = $true
$ws.Range("D28").Value = "https://www.railwayoperationsimulator.com/wp-content/uploads/2018/04/Eastbourne-Seaford-and-Newhaven-scaled.jpg"
$ws.Hyperlinks.Add($ws.Range("D28"), "https://www.railwayoperationsimulator.com/wp-content/uploads/2018/04/Eastbourne-Seaford-and-Newhaven-scaled.jpg") | Out-Null
$ws.Range("E28").Value = "-"

# --- Row 32: "Faversham to Ramsgate" ---
# Mark as added to map, record image URL as plain text, note "-"
$ws.Range("C32").Value = $true
$ws.Range("D32").Value = "https://www.railwayoperationsimulator.com/wp-content/uploads/2020/04/Faversham-to-Ramsgate.jpg"
$ws.Range("E32").Value = "-"

# --- Row 33: "Fenchurch Street" ---
# Add a note explaining the section is excluded
$ws.Range("E33").Value = "Not including this project"

# --- Row 36: "Grove Park, Hildenborough and Ashford" ---
# Was "WIP" - now finished, turn the existing link text into a live hyperlink
$ws.Range("C36").Value = $true
$ws.Hyperlinks.Add($ws.Range("D36"), "https://www.railwayoperationsimulator.com/catalog/community-projects/united-kingdom/grove-park-hildenborough-and-ashford") | Out-Null

# --- View state: active selection (scrolled so row 19 is the top of the view) ---
$ws.Application.Goto($ws.Range("A19"), $false) | Out-Null
$ws.Range("B32").Select() | Out-Null
